# Commit: "added new footer and a line"
#
# The diff appends two new paragraphs after the existing "Hi how are you"
# paragraph:
#   1. an empty paragraph (keeps the same en-IN language mark, no run)
#   2. a paragraph containing "I have added in new branch issue 53"
#
# Plain InsertParagraphAfter()/TypeParagraph() calls always materialize an
# empty <w:r> on the blank paragraph in this runtime, which the target
# markup does not have (it is a bare <w:pPr> with no run at all). To match
# the exact OOXML shape, build the two paragraphs as a WordprocessingML
# fragment and drop it in with Range.InsertXML at the end of the document -
# this is the same technique Range.InsertXML is designed for (inserting
# literal WordprocessingML), just expressed via the single-part package
# wrapper InsertXML expects.

$d = $word.ActiveDocument

$endRange = $d.Range($d.Content.End, $d.Content.End)

$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:lang w:val="en-IN"/>
              </w:rPr>
            </w:pPr>
          </w:p>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:lang w:val="en-IN"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:lang w:val="en-IN"/>
              </w:rPr>
              <w:t>I have added in new branch issue 53</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$endRange.InsertXML($xml)

Write-Output "Inserted blank paragraph + new line after the greeting."
